$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of an existing header cell (G1 - "sum") onto the new
# H1 header cell so the "Save" column header matches the other headers
# (bold font, border, centered alignment) using the same style record.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Set the new header label and its value for row 2
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 1
